# Update "paises" (countries) workbook: refresh COVID-19 stats snapshot.
# - Update the "last updated" timestamp
# - Refresh case counts for several countries (no ranking change)
# - Ecuador's case count overtakes Polonia/Luxemburgo/Filipinas -> re-rank rows 32-35
# - San Cristobal y Nieves overtakes Guinea-Bisau -> swap rows 174-175
# - Refresh Guyana's numbers (row 180)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title / timestamp cell
$ws.Range("A1").Value2 = "Datos actualizados a 31 de Marzo de 2020 a las 18:50"

# Estados Unidos (row 4) - same ranking, refreshed numbers
$ws.Cells.Item(4, 2).Value2 = 175669
$ws.Cells.Item(4, 3).Value2 = 11881
$ws.Cells.Item(4, 4).Value2 = 6241
$ws.Cells.Item(4, 5).Value2 = 166004
$ws.Cells.Item(4, 7).Value2 = 283
$ws.Cells.Item(4, 8).Value2 = 3424

# Brasil (row 21) - same ranking, refreshed numbers
$ws.Cells.Item(21, 2).Value2 = 4715
$ws.Cells.Item(21, 3).Value2 = 85
$ws.Cells.Item(21, 5).Value2 = 4420

# Noruega (row 22) - same ranking, refreshed numbers
$ws.Cells.Item(22, 5).Value2 = 4553
$ws.Cells.Item(22, 7).Value2 = 7
$ws.Cells.Item(22, 8).Value2 = 39

# Chequia (row 25) - same ranking, refreshed numbers
$ws.Cells.Item(25, 2).Value2 = 3138
$ws.Cells.Item(25, 3).Value2 = 137
$ws.Cells.Item(25, 5).Value2 = 3088

# Rumania (row 31) - same ranking, refreshed numbers
$ws.Cells.Item(31, 5).Value2 = 1946
$ws.Cells.Item(31, 7).Value2 = 14
$ws.Cells.Item(31, 8).Value2 = 79

# Rows 32-35 get re-ranked: Ecuador moves up (past Polonia, Luxemburgo, Filipinas)
# Row 32 becomes Ecuador with its refreshed numbers
$ws.Cells.Item(32, 1).Value2 = "Ecuador"
$ws.Cells.Item(32, 2).Value2 = 2240
$ws.Cells.Item(32, 3).Value2 = 274
$ws.Cells.Item(32, 4).Value2 = 54
$ws.Cells.Item(32, 5).Value2 = 2111
$ws.Cells.Item(32, 6).Value2 = 100
$ws.Cells.Item(32, 7).Value2 = 13
$ws.Cells.Item(32, 8).Value2 = 75

# Row 33 becomes Polonia (previous row 32 values, unchanged)
$ws.Cells.Item(33, 1).Value2 = "Polonia"
$ws.Cells.Item(33, 2).Value2 = 2215
$ws.Cells.Item(33, 3).Value2 = 160
$ws.Cells.Item(33, 4).Value2 = 7
$ws.Cells.Item(33, 5).Value2 = 2176
$ws.Cells.Item(33, 6).Value2 = 50
$ws.Cells.Item(33, 7).Value2 = 1
$ws.Cells.Item(33, 8).Value2 = 32

# Row 34 becomes Luxemburgo (previous row 33 values, unchanged)
$ws.Cells.Item(34, 1).Value2 = "Luxemburgo"
$ws.Cells.Item(34, 2).Value2 = 2178
$ws.Cells.Item(34, 3).Value2 = 190
$ws.Cells.Item(34, 4).Value2 = 80
$ws.Cells.Item(34, 5).Value2 = 2075
$ws.Cells.Item(34, 6).Value2 = 31
$ws.Cells.Item(34, 7).Value2 = 1
$ws.Cells.Item(34, 8).Value2 = 23

# Row 35 becomes Filipinas (previous row 34 values, unchanged)
$ws.Cells.Item(35, 1).Value2 = "Filipinas"
$ws.Cells.Item(35, 2).Value2 = 2084
$ws.Cells.Item(35, 3).Value2 = 538
$ws.Cells.Item(35, 4).Value2 = 49
$ws.Cells.Item(35, 5).Value2 = 1947
$ws.Cells.Item(35, 6).Value2 = 1
$ws.Cells.Item(35, 7).Value2 = 10
$ws.Cells.Item(35, 8).Value2 = 88

# Rows 174-175: San Cristobal y Nieves overtakes Guinea-Bisau
# Row 174 becomes San Cristobal y Nieves (previous row 175 values, unchanged)
$ws.Cells.Item(174, 1).Value2 = "San Cristobal y Nieves"
$ws.Cells.Item(174, 2).Value2 = 8
$ws.Cells.Item(174, 3).Value2 = 1
$ws.Cells.Item(174, 4).Value2 = 0
$ws.Cells.Item(174, 5).Value2 = 8
$ws.Cells.Item(174, 6).Value2 = 0
$ws.Cells.Item(174, 7).Value2 = 0
$ws.Cells.Item(174, 8).Value2 = 0

# Row 175 becomes Guinea-Bisau (previous row 174 values, unchanged)
$ws.Cells.Item(175, 1).Value2 = "Guinea-Bisau"
$ws.Cells.Item(175, 2).Value2 = 8
$ws.Cells.Item(175, 3).Value2 = 0
$ws.Cells.Item(175, 4).Value2 = 0
$ws.Cells.Item(175, 5).Value2 = 8
$ws.Cells.Item(175, 6).Value2 = 0
$ws.Cells.Item(175, 7).Value2 = 0
$ws.Cells.Item(175, 8).Value2 = 0

# Guyana (row 180) - same ranking, refreshed numbers
$ws.Cells.Item(180, 5).Value2 = 6
$ws.Cells.Item(180, 7).Value2 = 1
$ws.Cells.Item(180, 8).Value2 = 2
